$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 101; $row++) {
    $ws.Cells.Item($row, 5).Value = [double]($row - 2)
}
